$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for "Swedish" and "Uzbek" (search column A for an exact match)
$namesToRemove = @("Swedish", "Uzbek")
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

foreach ($name in $namesToRemove) {
    for ($r = $lastRow; $r -ge 1; $r--) {
        $cellValue = $ws.Cells.Item($r, 1).Value2
        if ($cellValue -eq $name) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
}

# Re-determine the data range (excluding header row 1) and sort by column B descending
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$sortRange = $ws.Range("A2:B$lastRow")
$keyRange = $ws.Range("B2:B$lastRow")

$sortRange.Sort($keyRange, 2)
